$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "VR Gaming experience" question column (H) and "Instructions time" column (I) ---
# Row 2 holds the variable definitions; row 3 holds the allowed value ranges ("Levels").
$ws.Range("H2").Value = "The participant's response to ""How much experience did you have with VR Gaming before today?"""
$ws.Range("I2").Value = "The total time it took the participant to read through the instructions"

$ws.Range("H3").Value = "0 - None at all \n 1 - I have used it, but not often \n 2 - I use it occasionally \n 3 - I use it often \n 4 - I use it all the time"

# --- Update the Levels definition for "Survey Version" (column C) ---
$ws.Range("C3").Value = "4, 5"

# --- Update the selected cell / scroll position on the sheet ---
[void]$ws.Range("I3").Select()

# --- Print options: center horizontally/vertically, show headings and gridlines when printing ---
$ws.PageSetup.CenterHorizontally = $true
$ws.PageSetup.CenterVertically = $true
$ws.PageSetup.PrintHeadings = $true
$ws.PageSetup.PrintGridlines = $true
